$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move K2 ("Pseudocode out Arthur Jump") to L2, preserving its style, then clear K2
$ws.Range("K2").Copy($ws.Range("L2"))
$ws.Range("K2").Clear()

# Move K3 ("Implement Arthur Jump") to L3, preserving its (default) style, then clear K3
$ws.Range("K3").Copy($ws.Range("L3"))
$ws.Range("K3").Clear()

# K5 switches fill style from white (s=1) to the yellow-highlight style used elsewhere (s=3),
# matching the formatting already used by e.g. I3/D3
$ws.Range("I3").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection to L2
$ws.Range("L2").Select()
